$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Jumpsuit,Caftan',
    'Parka,Jumpsuit',
    'Jumpsuit,Parka',
    'Parka,Gauchos',
    'Parka,Gauchos',
    'Jumpsuit,Halter',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Trunks',
    'Dress,Jumpsuit',
    'Jumpsuit,Kaftan',
    'Halter,Blazer',
    'Jumpsuit,Trunks',
    'Kaftan,Tee',
    'Parka,Jumpsuit',
    'Halter,Blazer',
    'Jumpsuit,Kaftan',
    'Halter,Jumpsuit',
    'Jumpsuit,Dress',
    'Halter,Blouse',
    'Blouse,Kaftan',
    'Blazer,Halter',
    'Blouse,Kaftan',
    'Blouse,Jumpsuit',
    'Jumpsuit,Dress',
    'Jumpsuit,Tee',
    'Jumpsuit,Kaftan',
    'Caftan,Blouse',
    'Jumpsuit,Halter',
    'Blazer,Blouse',
    'Jumpsuit,Blouse',
    'Jumpsuit,Halter',
    'Kaftan,Halter',
    'Jumpsuit,Blouse',
    'Cutoffs,Coat',
    'Jodhpurs,Trunks',
    'Halter,Caftan',
    'Jumpsuit,Kaftan',
    'Parka,Dress',
    'Halter,Jumpsuit',
    'Jumpsuit,Halter',
    'Caftan,Jumpsuit',
    'Halter,Jumpsuit',
    'Jodhpurs,Kaftan',
    'Blouse,Caftan',
    'Jumpsuit,Blouse',
    'Tee,Jumpsuit',
    'Jumpsuit,Blouse',
    'Halter,Blouse',
    'Caftan,Trunks',
    'Blouse,Jumpsuit',
    'Blazer,Halter',
    'Jumpsuit,Kaftan',
    'Caftan,Blazer',
    'Jumpsuit,Blazer',
    'Blouse,Caftan',
    'Jumpsuit,Trunks',
    'Parka,Blouse',
    'Jumpsuit,Blouse',
    'Jumpsuit,Caftan',
    'Jodhpurs,Trunks',
    'Caftan,Cutoffs',
    'Jumpsuit,Dress',
    'Halter,Trunks',
    'Caftan,Trunks',
    'Jumpsuit,Kaftan',
    'Tee,Kaftan',
    'Blazer,Halter',
    'Trunks,Jumpsuit',
    'Cutoffs,Trunks',
    'Parka,Blouse',
    'Parka,Blouse',
    'Jumpsuit,Tee',
    'Jumpsuit,Kaftan',
    'Sweatpants,Parka',
    'Trunks,Halter',
    'Halter,Jumpsuit',
    'Cutoffs,Turtleneck',
    'Blazer,Blouse',
    'Jumpsuit,Dress',
    'Jumpsuit,Kaftan',
    'Kaftan,Tee',
    'Kaftan,Tee',
    'Kaftan,Parka',
    'Jumpsuit,Top',
    'Blouse,Jumpsuit',
    'Kaftan,Jumpsuit'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
